# working on testng param to get testcase run while true
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from Sheet1 to TestData
$ws.Name = "TestData"

# Add a second test-case row
$ws.Range("A3").Value = "tc2"
$ws.Range("B3").Value = "false "
$ws.Range("C3").Value = "Minh "
$ws.Range("D3").Value = "Le "
$ws.Range("E3").Value = "Hoang"

# B2 currently holds the boolean TRUE; convert it to the literal text "true"
# so the TestNG runner reads it as a string parameter. Build it via a
# formula that yields the text, then paste-special as values so the cell
# ends up a plain shared-string cell (no cached formula, no quote-prefix).
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Formula = "=""true"""
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)

# Move the active selection to G7
$ws.Range("G7").Select()
